$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing")

# Clear the values of rows 3-7 (columns A-L), keeping their existing
# style/formatting intact - mirrors the pre-cleared rows below them.
$ws.Range("A3:L7").ClearContents()
